$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D (Price) to Text format so that numeric-looking
# strings (e.g. "2.00", "0.0825") are preserved exactly as text, matching
# the inlineStr cell type used in the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "50.044.21"
$ws.Range("E2").Value = "  +4.36%  "
$ws.Range("D3").Value = "2.663.69"
$ws.Range("E3").Value = "  +7.36%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "114.11"
$ws.Range("E5").Value = "  +8.65%  "
$ws.Range("D6").Value = "326.83"
$ws.Range("E6").Value = "  +3.33%  "
$ws.Range("E7").Value = "  +2.13%  "
$ws.Range("E9").Value = "  +4.02%  "
$ws.Range("D10").Value = "41.33"
$ws.Range("E10").Value = "  +6.64%  "
$ws.Range("D11").Value = "20.16"
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").Value = "0.0825"
$ws.Range("E12").Value = "  +3.34%  "
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").Value = "7.40"
$ws.Range("E14").Value = "  +4.69%  "
$ws.Range("D15").Value = "3.082.44"
$ws.Range("E15").Value = "  +7.31%  "
$ws.Range("D16").Value = "2.659.85"
$ws.Range("E16").Value = "  +7.07%  "
$ws.Range("D17").Value = "0.880"
$ws.Range("E17").Value = "  +6.68%  "
$ws.Range("D18").Value = "50.000.00"
$ws.Range("E18").Value = "  +4.42%  "
$ws.Range("D19").Value = "13.34"
$ws.Range("E19").Value = "  +5.39%  "
$ws.Range("E20").Value = "  +4.24%  "
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("E22").Value = "  +3.50%  "
$ws.Range("E23").Value = "  +2.79%  "
$ws.Range("D24").Value = "278.01"
$ws.Range("E24").Value = "  +2.66%  "
$ws.Range("E25").Value = "  +3.92%  "
$ws.Range("D26").Value = "26.93"
$ws.Range("E26").Value = "  +5.18%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "10.05"
$ws.Range("E28").Value = "  +3.99%  "
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("D30").Value = "36.75"
$ws.Range("E30").Value = "  +6.63%  "
$ws.Range("E31").Value = "  +3.65%  "
$ws.Range("D32").Value = "50.36"
$ws.Range("E32").Value = "  +2.09%  "
$ws.Range("E33").Value = "  +5.28%  "
$ws.Range("E34").Value = "  +4.48%  "
$ws.Range("E35").Value = "  +6.49%  "
$ws.Range("D36").Value = "5.10"
$ws.Range("E36").Value = "  +12.44%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  +8.19%  "
$ws.Range("D39").Value = "3.13"
$ws.Range("D40").Value = "125.67"
$ws.Range("E40").Value = "  +2.71%  "
$ws.Range("E41").Value = "  +2.53%  "
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").Value = "0.0319"
$ws.Range("E44").Value = "  +5.63%  "
$ws.Range("D45").Value = "2.111.86"
$ws.Range("E45").Value = "  +5.71%  "
$ws.Range("E46").Value = "  +6.52%  "
$ws.Range("E47").Value = "  +13.98%  "
$ws.Range("D48").Value = "2.00"
$ws.Range("E48").Value = "  +4.72%  "
$ws.Range("D49").Value = "9.12"
$ws.Range("E49").Value = "  +2.66%  "
$ws.Range("D50").Value = "5.37"
$ws.Range("E50").Value = "  +4.49%  "
$ws.Range("D51").Value = "59.69"
$ws.Range("E51").Value = "  +6.30%  "

# Restore the default (Normal) style on column D so no stray number-format
# style is left applied to the cells.
$ws.Range("D2:D51").Style = "Normal"
